$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Activate()

# Add a brand new row 3 of data (previously row3 only had formatted-but-empty E3/F3 cells).
# Cells are set in the same order that the shared-string table picks up new unique
# strings, so the resulting sharedStrings.xml ordering matches the target workbook.
$ws.Range("B3").Value = "extra corpus"
$ws.Range("O3").Value = "fixed capital issues with sentence"
$ws.Range("T3").Value = "top - like, options - die"
$ws.Range("U3").Value = "top - orphan, options - spiritual"

# Update W2: was "pass - stress" -> "fail - stress" (this introduces the "fail - stress" string)
$ws.Range("W2").Value = "fail - stress"

$ws.Range("V3").Value = "fail - morning"
$ws.Range("W3").Value = "fail - stress"
$ws.Range("X3").Value = "fail - look"
$ws.Range("Y3").Value = "matter"
$ws.Range("Z3").Value = "pass - top - case, options hand"
$ws.Range("AA3").Value = "pass - top"
$ws.Range("AB3").Value = "pass - top - basketball, options - outside"
$ws.Range("AC3").Value = "fail - best"

# Update view state to match the recorded selection after editing
$ws.Application.ActiveWindow.ScrollColumn = 18 # R is column 18
$ws.Range("AB4").Select()
